# Applies a 3-row cyclic rotation (row34 <- old row36, row35 <- old row34,
# row36 <- old row35) to the "Ryflodalen" observation rows of the sheet,
# matching the target diff (rows re-sorted by start time ascending:
# 12:46, 13:48, 14:37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 34: becomes the old row 36 data (Talltita / Poecile montanus) ----
$ws.Range("A34").Value = 130606943
$ws.Range("B34").Value = 58043
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 103021
$ws.Range("F34").Value = "Talltita"
$ws.Range("G34").Value = "Poecile montanus"
$ws.Range("H34").Value = "(Conrad von Baldenstein, 1827)"

$ws.Range("I34").Value = "'2"
$ws.Range("I34").Style = "Normal"

$ws.Range("L34").Value = "i par"
$ws.Range("M34").Value = "födosökande"

$ws.Range("Q34").Value = 485317
$ws.Range("R34").Value = 6939367
$ws.Range("Z34").Value = "12:46"
$ws.Range("AB34").Value = "12:46"

# ---- Row 35: becomes the old row 34 data (Vedticka / Fuscoporia viticola) ----
$ws.Range("A35").Value = 130608518
$ws.Range("B35").Value = 91771
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 5447
$ws.Range("F35").Value = "Vedticka"
$ws.Range("G35").Value = "Fuscoporia viticola"
$ws.Range("H35").Value = "(Schwein.) Murrill"

$ws.Range("Q35").Value = 485649
$ws.Range("R35").Value = 6939253
$ws.Range("Z35").Value = "13:48"
$ws.Range("AB35").Value = "13:48"

# ---- Row 36: becomes the old row 35 data (Skrovellav / Lobaria scrobiculata) ----
$ws.Range("A36").Value = 130609721
$ws.Range("B36").Value = 80349
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 2081
$ws.Range("F36").Value = "Skrovellav"
$ws.Range("G36").Value = "Lobaria scrobiculata"
$ws.Range("H36").Value = "(Scop.) DC."

# I36 no longer has an "Antal" value, and L36/M36 (Kon / Aktivitet) are cleared
$ws.Range("I36").Value = "'"
$ws.Range("I36").Style = "Normal"
$ws.Range("L36").ClearContents()
$ws.Range("M36").ClearContents()

$ws.Range("Q36").Value = 485397
$ws.Range("R36").Value = 6939386
$ws.Range("Z36").Value = "14:37"
$ws.Range("AB36").Value = "14:37"
